# Fix hardcoded mgt_rx / mgt_tx interface names (and related gty-suffixed
# identifiers) on the apex_ku15p_root_config sheet, matching the rename of
# the source .tab file to apex_ku15p_gty_root.tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("apex_ku15p_root_config")

$ws.Range("B2").Value  = "apex_ku15p_gty_root.tab"
$ws.Range("B7").Value  = "mgt_gty_rx"
$ws.Range("B8").Value  = "mgt_gty_tx"
$ws.Range("B9").Value  = "phalg_gty_tx"
$ws.Range("B10").Value = "drp_gty"
$ws.Range("B11").Value = "mgt_gty_interfaces.sv"
$ws.Range("B12").Value = "mgt_gty_module.sv"
$ws.Range("B13").Value = "common_gty_module.sv"
$ws.Range("B14").Value = "quad_gty_module.sv"
$ws.Range("B15").Value = "apex_ku15p_gty_serial_io"

# Update the saved cursor/selection position to match the new edit location.
$ws.Range("B11").Select()
